# Parandu 2025 workbook update - add newly uploaded flight rows (13-15)
# and a trailing marker cell (row 19) to the "Parandu" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: 12 Dec Flight2 3KM Beach -------------------------------------
$ws.Range("A13").Value = 46368
$ws.Range("A13").NumberFormat = "d-mmm"

$ws.Range("C13").Value = "cleaned_12 Dec_Flight2_3KM_Beach.csv"
$ws.Range("D13").Value = "12 Dec_Flight2_3KM.beach.bin-1239127.mat"

$ws.Range("E13").Value = 5960
$ws.Range("F13").Value = 2

# --- Row 14: 13 Dec Flight3 5KM up beach -----------------------------------
$ws.Range("A14").Value = 46369
$ws.Range("A14").NumberFormat = "d-mmm"

$ws.Range("B14").Value = 0.73263888888888884
$ws.Range("B14").NumberFormat = "h:mm"

$ws.Range("D14").Value = "13 Dec_Flight3_5KM_up_beach.bin-1863188.mat"
$ws.Range("C14").Value = "cleaned_13 Dec_Flight_5km_up_beach.csv"

$ws.Range("E14").Value = 1000
$ws.Range("F14").Value = 2

# --- Row 15: 14 Dec F2 afternoon 3KM beach ---------------------------------
$ws.Range("A15").Value = 46370
$ws.Range("A15").NumberFormat = "d-mmm"

$ws.Range("B15").Value = 0.53819444444444442
$ws.Range("B15").NumberFormat = "h:mm"

$ws.Range("C15").Value = "cleaned_14 Dec_F2_AFTERNOON_3KM_beach.csv"
$ws.Range("D15").Value = "2025-12-14 12-55-20.bin-1820656.mat"

$ws.Range("E15").Value = 1000
$ws.Range("F15").Value = 2

# --- Row 19: stray formatted cell left below the table ---------------------
$ws.Range("A19").NumberFormat = "d-mmm"

# Leave the selection where the author last left it.
$null = $ws.Range("D13").Select()
